$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 3; rows 3..36 (with their formatting/heights) shift
# down to become rows 4..37.
$ws.Rows.Item(3).Insert()

# Row 2 is updated in place with the newest draw result...
$ws.Range("B2").Value = "3 4 6 0" + [char]10 + "4 9 3 6" + [char]10 + "1 5 2 7" + [char]10 + "0 0 4 8"
$ws.Range("C2").Value = [char]0x2705 + " Direct: 14/3980 (0.35%)" + [char]10 + [char]0x2705 + " iBet: 14/215 (6.51%)"
# ...and its row height, which auto-expanded because of the wrapped text, is put
# back to the default (no explicit custom height), matching row 2's original state.
$ws.Rows.Item(2).AutoFit()

# The brand-new row 3 repeats the same date as row 2 and carries what used to be
# row 2's draw data.
$ws.Range("A3").Value = "28/6/2025 (Sat)"
$ws.Range("B3").Value = "2 0 5 7" + [char]10 + "6 2 8 8" + [char]10 + "9 1 3 0" + [char]10 + "7 6 1 4"
$ws.Range("C3").Value = [char]0x2705 + " Direct: 12/4302 (0.28%)" + [char]10 + [char]0x2705 + " iBet: 12/226 (5.31%)"
$ws.Range("B3:C3").WrapText = $true
$ws.Rows.Item(3).RowHeight = 60

# The data now naturally spans through row 8 (former row 7 shifted down), so the
# old blank placeholder row (now row 9) is removed to keep the trailing
# placeholder rows aligned as before.
$ws.Rows.Item(9).Delete()

# Row 11 gains a blank, styled C placeholder cell to match the rest of the table.
$ws.Range("C11").WrapText = $true

# A new blank placeholder row 37 is appended, mirroring row 36's B-only placeholder.
$ws.Range("B37").WrapText = $true
